$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# --- Data: fill in the 2020-04-26 row (row 92), which was previously blank ---
$ws.Cells.Item(92, 2).Value = 514
$ws.Cells.Item(92, 3).Value = 30028
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 6458

# --- Insert a new row for 2020-04-27 (row 93), pushing the footnote row down ---
# Insert() copies formatting down from the row above, matching row 92's styles.
$ws.Rows.Item(93).Insert()
$ws.Cells.Item(93, 1).Value = 43948

# --- Update the workbook-level print area to include the new row ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$96"
    }
}

# --- Update sheet view: zoom to 85% and move the selection ---
$ws.Activate()
$ws.Range("A92").Select()
$excel.ActiveWindow.Zoom = 85
